# Apply test-scenario updates to "Test Senaryoları" sheet (rows 3-7)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Senaryoları")

# Row 3: TC001 -> TC002 (Email Girişi)
$ws.Range("A3").Value = "TC002"
$ws.Range("B3").Value = "Email Girişi"
$ws.Range("C3").Value = "Email alanına geçerli email adresi girilmesi"
$ws.Range("F3").Value = 1

# Row 4: TC001 -> TC003 (Şifre Girişi)
$ws.Range("A4").Value = "TC003"
$ws.Range("B4").Value = "Şifre Girişi"
$ws.Range("C4").Value = "Şifre alanına geçerli şifre girilmesi"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "Şifre alanına tıkla"
$ws.Range("I4").Value = "password-input"
$ws.Range("J4").Value = "Tıkla"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = "Şifre alanı aktif"

# Row 5: TC001 -> TC004 (Giriş Butonu)
$ws.Range("A5").Value = "TC004"
$ws.Range("B5").Value = "Giriş Butonu"
$ws.Range("C5").Value = "Giriş butonuna tıklanarak giriş yapılması"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "Giriş butonuna tıkla"
$ws.Range("I5").Value = "login-btn"
$ws.Range("L5").Value = "Giriş yapıldı"

# Row 6: TC001 -> TC005 (Hatalı Email)
$ws.Range("A6").Value = "TC005"
$ws.Range("B6").Value = "Hatalı Email"
$ws.Range("C6").Value = "Geçersiz email ile giriş denemesi"
$ws.Range("D6").Value = "Orta"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "Email alanına geçersiz email gir"
$ws.Range("I6").Value = "email-input"
$ws.Range("K6").Value = "invalid@email"
$ws.Range("L6").Value = "Hata mesajı görüntülendi"

# Row 7: TC001 -> TC006 (Hatalı Şifre)
$ws.Range("A7").Value = "TC006"
$ws.Range("B7").Value = "Hatalı Şifre"
$ws.Range("C7").Value = "Geçersiz şifre ile giriş denemesi"
$ws.Range("D7").Value = "Orta"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "Şifre alanına geçersiz şifre gir"
$ws.Range("I7").Value = "password-input"
$ws.Range("J7").Value = "Yaz"
$ws.Range("K7").Value = "wrongpass"
$ws.Range("L7").Value = "Hata mesajı görüntülendi"
